$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.570.67"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "2.051.55"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'245.45"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("D6").Value = "'0.659"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'54.40"
$ws.Range("E8").Value = "  -7.95%  "
$ws.Range("D9").Value = "'60.43"
$ws.Range("E9").Value = "  +1.41%  "
$ws.Range("D10").Value = "'0.364"
$ws.Range("E10").Value = "  -3.38%  "
$ws.Range("D11").Value = "'0.0749"
$ws.Range("E11").Value = "  -4.22%  "
$ws.Range("E12").Value = "  -3.61%  "
$ws.Range("E13").Value = "  +8.39%  "
$ws.Range("D14").Value = "'14.72"
$ws.Range("E14").Value = "  -4.81%  "
$ws.Range("D15").Value = "2.352.78"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").Value = "'5.45"
$ws.Range("E16").Value = "  -4.85%  "
$ws.Range("D17").Value = "2.067.34"
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("D18").Value = "36.470.28"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("D19").Value = "'17.06"
$ws.Range("E19").Value = "  -6.49%  "
$ws.Range("D20").Value = "'71.83"
$ws.Range("E20").Value = "  -2.90%  "
$ws.Range("D21").Value = "0.0₃0857"
$ws.Range("E21").Value = "  -4.43%  "
$ws.Range("D22").Value = "'237.90"
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("D23").Value = "'5.22"
$ws.Range("E23").Value = "  -4.46%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").Value = "'2.38"
$ws.Range("E25").Value = "  -2.78%  "
$ws.Range("D26").Value = "'2.25"
$ws.Range("E26").Value = "  +4.19%  "
$ws.Range("D27").Value = "'164.98"
$ws.Range("E27").Value = "  -3.21%  "
$ws.Range("D28").Value = "'9.17"
$ws.Range("E28").Value = "  -10.51%  "
$ws.Range("D29").Value = "'20.03"
$ws.Range("E29").Value = "  -0.59%  "
$ws.Range("E30").Value = "  -2.63%  "
$ws.Range("D31").Value = "'1.22"
$ws.Range("E31").Value = "  +9.06%  "
$ws.Range("E32").Value = "  -8.54%  "
$ws.Range("D33").Value = "'4.46"
$ws.Range("E33").Value = "  -5.57%  "
$ws.Range("D34").Value = "'0.0592"
$ws.Range("E34").Value = "  -4.45%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "'0.0865"
$ws.Range("E36").Value = "  +2.63%  "
$ws.Range("D37").Value = "'1.82"
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("D38").Value = "'2.22"
$ws.Range("E38").Value = "  -6.45%  "
$ws.Range("B39").Value = "THORChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D39").Value = "'5.01"
$ws.Range("E39").Value = "  -4.58%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.24"
$ws.Range("E40").Value = "  -7.45%  "
$ws.Range("E41").Value = "  -5.65%  "
$ws.Range("D42").Value = "'0.0214"
$ws.Range("E42").Value = "  -4.60%  "
$ws.Range("E43").Value = "  -5.33%  "
$ws.Range("D44").Value = "'94.13"
$ws.Range("E44").Value = "  -3.85%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.412.40"
$ws.Range("E45").Value = "  +8.26%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "'0.0906"
$ws.Range("E46").Value = "  -5.45%  "
$ws.Range("D47").Value = "'15.86"
$ws.Range("E47").Value = "  -6.82%  "
$ws.Range("D48").Value = "'7.45"
$ws.Range("E48").Value = "  +9.21%  "
$ws.Range("E49").Value = "  +1.49%  "
$ws.Range("D50").Value = "'2.25"
$ws.Range("E50").Value = "  -5.24%  "
$ws.Range("D51").Value = "2.233.32"
$ws.Range("E51").Value = "  -0.58%  "
